$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.418.26'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.390.70'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.20%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.93'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.21'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.76%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.538'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.99%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.392.24'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.76%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.24%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.04%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.338'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.70'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.847.10'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.89%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.563.39'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.80%  '

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.398.76'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.76%  '

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.32'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +15.07%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '325.18'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.85%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.76%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.79'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -7.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.50'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '555.31'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.01'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -11.99%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0910'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.91'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.15%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.50%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.71%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.34%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.41%  '

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '153.75'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.72%  '

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.41'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.367'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.51'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.16%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.06'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.30'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.51%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.79%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '142.59'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.48'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.586'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0499'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.90'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.69%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.37%  '
